$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-05 16:18:05'
$ws.Range("H2").Value = "'94%"
$ws.Range("K2").Value = '3.4 MJ/m2'
$ws.Range("E3").Value = '2026-02-05 16:18:08'
$ws.Range("K3").Value = '5.9 MJ/m2'
$ws.Range("O3").Value = '-2.1 °C'
$ws.Range("E4").Value = '2026-02-05 16:18:10'
$ws.Range("J4").Value = '991.0 hPa'
$ws.Range("K4").Value = '6.5 MJ/m2'
$ws.Range("L4").Value = '61.2 km/h - 288º 15:40 TU'
$ws.Range("O4").Value = '10.5 °C'
$ws.Range("E5").Value = '2026-02-05 16:18:12'
$ws.Range("H5").Value = "'75%"
$ws.Range("J5").Value = '991.0 hPa'
$ws.Range("K5").Value = '6.4 MJ/m2'
$ws.Range("O5").Value = '8.8 °C'
$ws.Range("E6").Value = '2026-02-05 16:18:15'
$ws.Range("H6").Value = "'74%"
$ws.Range("J6").Value = '992.4 hPa'
$ws.Range("K6").Value = '3.5 MJ/m2'
$ws.Range("L6").Value = '31.3 km/h - 251º 15:46 TU'
$ws.Range("M6").Value = '15.5 °C 15:59 TU'
$ws.Range("O6").Value = '12.3 °C'
$ws.Range("E7").Value = '2026-02-05 16:18:18'
$ws.Range("J7").Value = '992.3 hPa'
$ws.Range("K7").Value = '2.7 MJ/m2'
$ws.Range("M7").Value = '13.1 °C 15:59 TU'
$ws.Range("O7").Value = '10.0 °C'
$ws.Range("E8").Value = '2026-02-05 16:18:21'
$ws.Range("K8").Value = '5.4 MJ/m2'
$ws.Range("O8").Value = '7.6 °C'
$ws.Range("E9").Value = '2026-02-05 16:18:23'
$ws.Range("O9").Value = '1.7 °C'
$ws.Range("E10").Value = '2026-02-05 16:18:25'
$ws.Range("H10").Value = "'92%"
$ws.Range("O10").Value = '6.7 °C'
$ws.Range("E11").Value = '2026-02-05 16:18:28'
$ws.Range("I11").Value = '8.2 mm'
$ws.Range("J11").Value = '995.9 hPa'
$ws.Range("K11").Value = '3.0 MJ/m2'
$ws.Range("E12").Value = '2026-02-05 16:18:31'
$ws.Range("K12").Value = '3.6 MJ/m2'
$ws.Range("M12").Value = '15.5 °C 15:59 TU'
$ws.Range("O12").Value = '8.8 °C'
$ws.Range("E13").Value = '2026-02-05 16:18:34'
$ws.Range("H13").Value = "'85%"
$ws.Range("O13").Value = '7.0 °C'
$ws.Range("E14").Value = '2026-02-05 16:18:36'
$ws.Range("H14").Value = "'68%"
$ws.Range("I14").Value = '5.1 mm'
$ws.Range("K14").Value = '1.6 MJ/m2'
$ws.Range("E15").Value = '2026-02-05 16:18:39'
$ws.Range("J15").Value = '991.5 hPa'
$ws.Range("K15").Value = '6.4 MJ/m2'
$ws.Range("O15").Value = '6.6 °C'
$ws.Range("E16").Value = '2026-02-05 16:18:42'
$ws.Range("M16").Value = '5.9 °C 15:59 TU'
$ws.Range("O16").Value = '3.3 °C'
$ws.Range("E17").Value = '2026-02-05 16:18:45'
$ws.Range("J17").Value = '995.8 hPa'
$ws.Range("K17").Value = '1.9 MJ/m2'
$ws.Range("O17").Value = '0.7 °C'
$ws.Range("E18").Value = '2026-02-05 16:18:47'
$ws.Range("K18").Value = '1.2 MJ/m2'
$ws.Range("O18").Value = '-4.5 °C'
$ws.Range("E19").Value = '2026-02-05 16:18:50'
$ws.Range("J19").Value = '992.8 hPa'
$ws.Range("K19").Value = '4.7 MJ/m2'
$ws.Range("O19").Value = '6.9 °C'
$ws.Range("E20").Value = '2026-02-05 16:18:53'
$ws.Range("H20").Value = "'71%"
$ws.Range("I20").Value = '4.0 mm'
$ws.Range("K20").Value = '1.0 MJ/m2'
$ws.Range("O20").Value = '-1.7 °C'
$ws.Range("E21").Value = '2026-02-05 16:18:56'
$ws.Range("H21").Value = "'86%"
$ws.Range("J21").Value = '991.9 hPa'
$ws.Range("K21").Value = '5.7 MJ/m2'
$ws.Range("M21").Value = '13.3 °C 15:59 TU'
$ws.Range("O21").Value = '5.2 °C'
$ws.Range("E22").Value = '2026-02-05 16:18:59'
$ws.Range("K22").Value = '5.1 MJ/m2'
$ws.Range("O22").Value = '7.4 °C'
$ws.Range("E23").Value = '2026-02-05 16:19:02'
$ws.Range("H23").Value = "'87%"
$ws.Range("J23").Value = '991.0 hPa'
$ws.Range("K23").Value = '3.3 MJ/m2'
$ws.Range("L23").Value = '34.2 km/h - 288º 15:46 TU'
$ws.Range("M23").Value = '15.2 °C 15:57 TU'
$ws.Range("O23").Value = '7.6 °C'
$ws.Range("E24").Value = '2026-02-05 16:19:05'
$ws.Range("H24").Value = "'78%"
$ws.Range("J24").Value = '990.0 hPa'
$ws.Range("K24").Value = '5.1 MJ/m2'
$ws.Range("O24").Value = '9.9 °C'
$ws.Range("E25").Value = '2026-02-05 16:19:07'
$ws.Range("J25").Value = '994.9 hPa'
$ws.Range("K25").Value = '3.1 MJ/m2'
$ws.Range("E26").Value = '2026-02-05 16:19:10'
$ws.Range("K26").Value = '4.1 MJ/m2'
$ws.Range("L26").Value = '59.8 km/h - 16º 15:43 TU'
$ws.Range("E27").Value = '2026-02-05 16:19:13'
$ws.Range("H27").Value = "'90%"
$ws.Range("J27").Value = '991.3 hPa'
$ws.Range("K27").Value = '3.3 MJ/m2'
$ws.Range("M27").Value = '15.0 °C 15:37 TU'
$ws.Range("O27").Value = '7.8 °C'
$ws.Range("E28").Value = '2026-02-05 16:19:16'
$ws.Range("H28").Value = "'96%"
$ws.Range("J28").Value = '994.3 hPa'
$ws.Range("M28").Value = '7.0 °C 15:57 TU'
$ws.Range("O28").Value = '1.6 °C'
$ws.Range("E29").Value = '2026-02-05 16:19:19'
$ws.Range("K29").Value = '4.3 MJ/m2'
$ws.Range("M29").Value = '14.5 °C 15:57 TU'
$ws.Range("O29").Value = '7.8 °C'
$ws.Range("E30").Value = '2026-02-05 16:19:21'
$ws.Range("H30").Value = "'66%"
$ws.Range("I30").Value = '4.8 mm'
$ws.Range("K30").Value = '1.2 MJ/m2'
$ws.Range("O30").Value = '-2.3 °C'
$ws.Range("E31").Value = '2026-02-05 16:19:24'
$ws.Range("I31").Value = '17.1 mm'
$ws.Range("J31").Value = '994.9 hPa'
$ws.Range("E32").Value = '2026-02-05 16:19:27'
$ws.Range("H32").Value = "'85%"
$ws.Range("J32").Value = '992.3 hPa'
$ws.Range("K32").Value = '5.8 MJ/m2'
$ws.Range("M32").Value = '18.7 °C 15:32 TU'
$ws.Range("O32").Value = '11.3 °C'
$ws.Range("E33").Value = '2026-02-05 16:19:29'
$ws.Range("H33").Value = "'90%"
$ws.Range("O33").Value = '7.9 °C'
$ws.Range("E34").Value = '2026-02-05 16:19:31'
$ws.Range("K34").Value = '2.5 MJ/m2'
$ws.Range("M34").Value = '9.1 °C 15:54 TU'
$ws.Range("O34").Value = '2.6 °C'
$ws.Range("E35").Value = '2026-02-05 16:19:34'
$ws.Range("I35").Value = '3.3 mm'
$ws.Range("K35").Value = '2.7 MJ/m2'
$ws.Range("M35").Value = '-1.3 °C 15:59 TU'
$ws.Range("E36").Value = '2026-02-05 16:19:37'
$ws.Range("H36").Value = "'87%"
$ws.Range("J36").Value = '992.4 hPa'
$ws.Range("K36").Value = '9.5 MJ/m2'
$ws.Range("O36").Value = '9.6 °C'
